$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Price" (column D) values. The $true flag marks values that
# look numeric (e.g. "1.001") and therefore need the cell pre-formatted
# as Text so Excel stores them as the literal string rather than a number
# (values such as "27.219.03" already fail numeric parsing and need no help).
$priceUpdates = @(
    @{Row=2; Value='27.219.03'; NumericLooking=$false}
    @{Row=3; Value='1.905.57'; NumericLooking=$false}
    @{Row=4; Value='1.001'; NumericLooking=$true}
    @{Row=5; Value='305.93'; NumericLooking=$true}
    @{Row=6; Value='0.9996'; NumericLooking=$true}
    @{Row=7; Value='0.5389'; NumericLooking=$true}
    @{Row=8; Value='0.3804'; NumericLooking=$true}
    @{Row=9; Value='0.07288'; NumericLooking=$true}
    @{Row=10; Value='22.28'; NumericLooking=$true}
    @{Row=11; Value='0.9045'; NumericLooking=$true}
    @{Row=12; Value='0.08188'; NumericLooking=$true}
    @{Row=13; Value='95.72'; NumericLooking=$true}
    @{Row=14; Value='5.336'; NumericLooking=$true}
    @{Row=15; Value='1.003'; NumericLooking=$true}
    @{Row=18; Value='0.9999'; NumericLooking=$true}
    @{Row=19; Value='27.251.82'; NumericLooking=$false}
    @{Row=20; Value='1.130.14'; NumericLooking=$false}
    @{Row=21; Value='5.045'; NumericLooking=$true}
    @{Row=23; Value='6.525'; NumericLooking=$true}
    @{Row=24; Value='148.57'; NumericLooking=$true}
    @{Row=25; Value='2.306'; NumericLooking=$true}
    @{Row=27; Value='1.752'; NumericLooking=$true}
    @{Row=28; Value='116.57'; NumericLooking=$true}
    @{Row=29; Value='4.848'; NumericLooking=$true}
    @{Row=30; Value='4.713'; NumericLooking=$true}
    @{Row=31; Value='0.09216'; NumericLooking=$true}
    @{Row=32; Value='0.8310'; NumericLooking=$true}
    @{Row=33; Value='0.05075'; NumericLooking=$true}
    @{Row=35; Value='3.006'; NumericLooking=$true}
    @{Row=36; Value='3.326'; NumericLooking=$true}
    @{Row=37; Value='2.661'; NumericLooking=$true}
    @{Row=38; Value='0.5921'; NumericLooking=$true}
    @{Row=39; Value='0.02001'; NumericLooking=$true}
    @{Row=40; Value='1.083'; NumericLooking=$true}
    @{Row=41; Value='9.298'; NumericLooking=$true}
    @{Row=42; Value='6.656'; NumericLooking=$true}
    @{Row=43; Value='116.57'; NumericLooking=$true}
    @{Row=44; Value='0.5121'; NumericLooking=$true}
    @{Row=47; Value='0.9995'; NumericLooking=$true}
    @{Row=48; Value='1.641'; NumericLooking=$true}
    @{Row=49; Value='38.31'; NumericLooking=$true}
    @{Row=50; Value='0.06124'; NumericLooking=$true}
    @{Row=51; Value='63.43'; NumericLooking=$true}
)

foreach ($u in $priceUpdates) {
    $cell = $ws.Range("D" + $u.Row)
    if ($u.NumericLooking) {
        $cell.NumberFormat = '@'
        $cell.Value = $u.Value
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $u.Value
    }
}

# Updated "Volume(1h)" (column E) values - plain text, never numeric-looking
# (they keep the padding spaces and "%" sign), so no text-forcing is needed.
$volumeUpdates = @(
    @{Row=2; Value='  +0.31%  '}
    @{Row=3; Value='  +0.79%  '}
    @{Row=4; Value='  -0.16%  '}
    @{Row=5; Value='  -0.29%  '}
    @{Row=6; Value='  -0.17%  '}
    @{Row=7; Value='  +3.41%  '}
    @{Row=8; Value='  +1.39%  '}
    @{Row=9; Value='  +0.43%  '}
    @{Row=10; Value='  +5.84%  '}
    @{Row=11; Value='  +0.79%  '}
    @{Row=12; Value='  +0.29%  '}
    @{Row=13; Value='  -0.54%  '}
    @{Row=14; Value='  +1.33%  '}
    @{Row=15; Value='  +0.03%  '}
    @{Row=16; Value='  +2.42%  '}
    @{Row=17; Value='  +0.99%  '}
    @{Row=18; Value='  -0.16%  '}
    @{Row=19; Value='  +0.28%  '}
    @{Row=20; Value='  -40.27%  '}
    @{Row=21; Value='  -0.66%  '}
    @{Row=22; Value='  +1.03%  '}
    @{Row=23; Value='  +2.09%  '}
    @{Row=24; Value='  +0.78%  '}
    @{Row=25; Value='  +0.37%  '}
    @{Row=26; Value='  +1.24%  '}
    @{Row=27; Value='  +1.05%  '}
    @{Row=28; Value='  +1.47%  '}
    @{Row=29; Value='  +1.48%  '}
    @{Row=30; Value='  -3.70%  '}
    @{Row=31; Value='  -0.07%  '}
    @{Row=32; Value='  +5.20%  '}
    @{Row=33; Value='  +0.70%  '}
    @{Row=34; Value='  +0.78%  '}
    @{Row=35; Value='  +1.21%  '}
    @{Row=36; Value='  -2.87%  '}
    @{Row=37; Value='  +3.42%  '}
    @{Row=38; Value='  +4.78%  '}
    @{Row=39; Value='  +0.77%  '}
    @{Row=40; Value='  +0.93%  '}
    @{Row=41; Value='  +3.67%  '}
    @{Row=42; Value='  +1.81%  '}
    @{Row=43; Value='  +0.63%  '}
    @{Row=44; Value='  +5.68%  '}
    @{Row=45; Value='  +1.00%  '}
    @{Row=46; Value='  +2.03%  '}
    @{Row=47; Value='  -0.14%  '}
    @{Row=48; Value='  +1.68%  '}
    @{Row=50; Value='  +3.18%  '}
)

foreach ($u in $volumeUpdates) {
    $ws.Range("E" + $u.Row).Value = $u.Value
}
